$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (new TPM-derived stats)
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.02550355756355556
$ws.Range("R2").Value = 0.229532018072
$ws.Range("S2").Value = 0.3931387525216601
$ws.Range("T2").Value = 0.39313875252166

# Update row 3 (new TPM-derived specificity values)
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("S3").Value = 0.60686124747834
$ws.Range("T3").Value = 0.60686124747834

# Row 4 (Kng1/Itgb2 -> MuSCs) no longer exists with new data; remove it
$ws.Rows.Item(4).Delete()
